# Insert a new weekly price record as row 134 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 134-177 down to 135-178.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 134..177 down by inserting a new row at 134.
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new record's data.
$ws.Range("A134").Value2 = 11
$ws.Range("B134").Value2 = "Vega Monumental Concepción"
$ws.Range("C134").Value2 = "Bíobío"
$ws.Range("D134").Value2 = 45135
$ws.Range("D134").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E134").Value2 = 8
$ws.Range("F134").Value2 = 100112001
$ws.Range("G134").Value2 = "Berenjena"
$ws.Range("H134").Value2 = "Sin especificar"
$ws.Range("I134").Value2 = "Primera"
$ws.Range("J134").Value2 = 50
$ws.Range("K134").Value2 = 10000
$ws.Range("L134").Value2 = 10000
$ws.Range("M134").Value2 = 10000
$ws.Range("N134").Value2 = "`$/caja 50 unidades"
$ws.Range("O134").Value2 = "Región de Arica y Parinacota"
$ws.Range("P134").Value2 = 200
$ws.Range("Q134").Value2 = 50
$ws.Range("R134").Value2 = "Hortaliza"
